$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths: col G (7) 4.140625 -> 5.140625, col L (12) 3.140625 -> 4.140625 (raw OOXML "width" units).
# The ColumnWidth COM property here quantizes internally to 1/6-character steps
# (stored_width = round(ColumnWidth*6)/6 + 5/6), so these inputs land on the nearest
# reachable bucket (5.16666... / 4.16666...), the closest possible approximation of the target widths.
$ws.Columns.Item(7).ColumnWidth = 4.333333333333333
$ws.Columns.Item(12).ColumnWidth = 3.3333333333333335

# Update data values in columns G-L for rows 1-17
$ws.Range("G1").Value = 221
$ws.Range("H1").Value = 26.8661671
$ws.Range("I1").Value = 0.0000616487089080131
$ws.Range("J1").Value = 0.0000007648595993813746
$ws.Range("K1").Value = 0
$ws.Range("L1").Value = 0
$ws.Range("G2").Value = 269
$ws.Range("H2").Value = 73.478886
$ws.Range("I2").Value = 0.00006551989470948705
$ws.Range("J2").Value = 0.0000005057287378485495
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("G3").Value = 269
$ws.Range("H3").Value = 82.1382684
$ws.Range("I3").Value = 0.00006676530033189465
$ws.Range("J3").Value = 0.0000005170180650411631
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("G4").Value = 275
$ws.Range("H4").Value = 143.1748107
$ws.Range("I4").Value = 0.00006060388280171303
$ws.Range("J4").Value = 0.0000006556261854759074
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("G5").Value = 238
$ws.Range("H5").Value = 79.3551348
$ws.Range("I5").Value = 0.00006771588698573083
$ws.Range("J5").Value = 0.0000006178503110070885
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 57
$ws.Range("G6").Value = 779
$ws.Range("H6").Value = 404.4532273
$ws.Range("I6").Value = 0.00006515203979939876
$ws.Range("J6").Value = 0.00000018010085309868908
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("G7").Value = 1262
$ws.Range("H7").Value = 280.9737301
$ws.Range("I7").Value = 0.0001559713646286376
$ws.Range("J7").Value = -0.00000164414825912559
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("G8").Value = 285
$ws.Range("H8").Value = 87.028548
$ws.Range("I8").Value = 0.0000679872891180544
$ws.Range("J8").Value = 0.0000004336788694957792
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 61
$ws.Range("G9").Value = 312
$ws.Range("H9").Value = 94.7836417
$ws.Range("I9").Value = 0.0000681160423607885
$ws.Range("J9").Value = 0.000000606933726764572
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 83
$ws.Range("G10").Value = 369
$ws.Range("H10").Value = 176.2547574
$ws.Range("I10").Value = 0.00018531213506078892
$ws.Range("J10").Value = -0.000004398869262743244
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 44
$ws.Range("G11").Value = 269
$ws.Range("H11").Value = 161.8426977
$ws.Range("I11").Value = 0.000784156263761826
$ws.Range("J11").Value = -0.00003060860791280304
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 172
$ws.Range("G12").Value = 269
$ws.Range("H12").Value = 243.2040577
$ws.Range("I12").Value = 0.0003568319962006239
$ws.Range("J12").Value = -0.00003976692669122803
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 172
$ws.Range("G13").Value = 269
$ws.Range("H13").Value = 236.4909864
$ws.Range("I13").Value = 0.001952341701944027
$ws.Range("J13").Value = -0.000030092198418340354
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 172
$ws.Range("G14").Value = 209
$ws.Range("H14").Value = 10.7656706
$ws.Range("I14").Value = 0.00006345625439818647
$ws.Range("J14").Value = 0.0000007054478258456188
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("G15").Value = 269
$ws.Range("H15").Value = 115.4354966
$ws.Range("I15").Value = 0.00006676530033189465
$ws.Range("J15").Value = 0.0000005170180650411631
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("G16").Value = 251
$ws.Range("H16").Value = 92.0375708
$ws.Range("I16").Value = 0.00020161746254965252
$ws.Range("J16").Value = -0.000023781491897126006
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 83
$ws.Range("G17").Value = 269
$ws.Range("H17").Value = 221.2334096
$ws.Range("I17").Value = 0.000784156263761826
$ws.Range("J17").Value = -0.00003060860791280304
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 172
